$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells (row 1) ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Update column C (GDP) values for data rows ---
$ws.Range("C2").Value = 7854.952374701078
$ws.Range("C3").Value = 2934.187009790061
$ws.Range("C4").Value = 2870.311589353206
$ws.Range("C5").Value = 1873.394108966653
$ws.Range("C6").Value = 1460.056109840828
$ws.Range("C7").Value = 5191.140356354663
$ws.Range("C8").Value = 1909.084588129339
$ws.Range("C9").Value = 10594.98659239237
$ws.Range("C10").Value = 6128.19547247793
$ws.Range("C11").Value = 4547.50930098406
$ws.Range("C12").Value = 4729.735976516416
$ws.Range("C13").Value = 749.9184730334429
$ws.Range("C14").Value = 1250.795760575873
$ws.Range("C15").Value = 3587.883798243964
$ws.Range("C16").Value = 471.181692645893
$ws.Range("C17").Value = 488.421401781569
$ws.Range("C18").Value = 1036.533951644687
$ws.Range("C19").Value = 647.8358464534491
$ws.Range("C20").Value = 534.0443668101981
$ws.Range("C21").Value = 341.5541149051794
$ws.Range("C22").Value = 1401.47747416771
$ws.Range("C23").Value = 2898.942214704482
$ws.Range("C24").Value = 665.6274194933962
$ws.Range("C25").Value = 1503.870423231357
$ws.Range("C26").Value = 5555.389721901988
$ws.Range("C27").Value = 1955.461557360978
$ws.Range("C28").Value = 5082.354756663512
$ws.Range("C29").Value = 2217.474008566157
$ws.Range("C30").Value = 1317.890706178356
$ws.Range("C31").Value = 3579.960081455846
$ws.Range("C32").Value = 17610.30663334184
$ws.Range("C33").Value = 505.2384587280311
$ws.Range("C34").Value = 1037.747039954749
$ws.Range("C35").Value = 1132.548400540401
$ws.Range("C36").Value = 694.6031345426339
$ws.Range("C37").Value = 711.1128122770988
$ws.Range("C38").Value = 553.2014555484933
$ws.Range("C39").Value = 369.2024078290272
$ws.Range("C40").Value = 1591.56825353313
$ws.Range("C41").Value = 3083.80337578809
$ws.Range("C42").Value = 691.8942672110555
$ws.Range("C43").Value = 1577.487171555845
$ws.Range("C44").Value = 5660.517066940175
$ws.Range("C45").Value = 2024.117324382548
$ws.Range("C46").Value = 11627.81065059172
$ws.Range("C47").Value = 2264.394087033834
$ws.Range("C48").Value = 838.188712186125
$ws.Range("C49").Value = 1657.651524528445
$ws.Range("C50").Value = 5745.422744292303
$ws.Range("C51").Value = 701.4459636783288
$ws.Range("C52").Value = 2094.024217383061
$ws.Range("C53").Value = 2999.422762626143
$ws.Range("C54").Value = 6911.59200404802
$ws.Range("C55").Value = 3748.449444923865
$ws.Range("C56").Value = 513.7390871590731
$ws.Range("C57").Value = 2860.874335573629
$ws.Range("C58").Value = 1057.667740311969
$ws.Range("C59").Value = 719.6981727039259
$ws.Range("C60").Value = 684.6474015015979
$ws.Range("C61").Value = 1223.631935023299
$ws.Range("C62").Value = 573.9239887389259
$ws.Range("C63").Value = 389.9389667216314
$ws.Range("C64").Value = 1745.10167474004
$ws.Range("C65").Value = 2379.668184479739
$ws.Range("C66").Value = 1716.389195271215
$ws.Range("C67").Value = 5955.175904294275
$ws.Range("C68").Value = 720.7128711178943
$ws.Range("C69").Value = 2201.396847776877
$ws.Range("C70").Value = 3056.152683606517
$ws.Range("C71").Value = 7200.731056811853
$ws.Range("C72").Value = 3796.882621798447
$ws.Range("C73").Value = 534.5063430177229
$ws.Range("C74").Value = 2887.250212489506
$ws.Range("C75").Value = 1102.527430026863
$ws.Range("C76").Value = 731.5588677998553
$ws.Range("C77").Value = 680.3923729568069
$ws.Range("C78").Value = 1299.811672673934
$ws.Range("C79").Value = 593.1620921048029
$ws.Range("C80").Value = 419.1838602515346
$ws.Range("C81").Value = 1778.60982580794
$ws.Range("C82").Value = 2497.68592515536
$ws.Range("C83").Value = 3843.198240901342
$ws.Range("C84").Value = 1140.447753778042
$ws.Range("C85").Value = 2286.013198234259
$ws.Range("C86").Value = 11951.20944634967
$ws.Range("C87").Value = 7449.08671983612
$ws.Range("C88").Value = 3008.669179463094
$ws.Range("C89").Value = 1379.14068216006
$ws.Range("C90").Value = 5412.131646018807
$ws.Range("C91").Value = 449.4203771491282
$ws.Range("C92").Value = 1627.760281433693
$ws.Range("C93").Value = 3137.260298393558
$ws.Range("C94").Value = 730.3063521039821
$ws.Range("C95").Value = 707.8672001573369
$ws.Range("C96").Value = 711.3043470146426
$ws.Range("C97").Value = 1775.027517189621
$ws.Range("C98").Value = 6301.696269820412
$ws.Range("C99").Value = 567.8342670439314
$ws.Range("C101").Value = 612.1489724037899
$ws.Range("C102").Value = 3748.320622951519
$ws.Range("C103").Value = 1128.996380299766
$ws.Range("C104").Value = 2361.056581219794
$ws.Range("C105").Value = 11431.15448084494
$ws.Range("C106").Value = 7580.275568826287
$ws.Range("C107").Value = 3012.536723186288
$ws.Range("C108").Value = 1463.71052702022
$ws.Range("C109").Value = 5330.539154475424
$ws.Range("C110").Value = 482.6390663355013
$ws.Range("C111").Value = 1625.905825842452
$ws.Range("C112").Value = 3210.869677115934
$ws.Range("C113").Value = 729.1196658666737
$ws.Range("C114").Value = 729.7808175407341
$ws.Range("C115").Value = 731.9993357350996
$ws.Range("C116").Value = 1836.014008604312
$ws.Range("C117").Value = 6661.86504232374
$ws.Range("C118").Value = 441.1376640642927
$ws.Range("C120").Value = 630.9372503341563
$ws.Range("C121").Value = 3530.309422482455
$ws.Range("C122").Value = 1134.924536209078
$ws.Range("C123").Value = 2425.561644739583
$ws.Range("C124").Value = 10965.97426143915
$ws.Range("C125").Value = 7633.969039669125
$ws.Range("C126").Value = 2854.757682901436
$ws.Range("C127").Value = 1529.507453727912
$ws.Range("C128").Value = 5176.058803160127
$ws.Range("C129").Value = 514.0573067519859
$ws.Range("C130").Value = 1644.598009122967
$ws.Range("C131").Value = 3242.636921959078
$ws.Range("C132").Value = 729.8559996981501
$ws.Range("C133").Value = 749.2194349876407
$ws.Range("C134").Value = 10205.79575322194
$ws.Range("C135").Value = 729.6614300490079
$ws.Range("C136").Value = 1895.214690888655
$ws.Range("C137").Value = 359.6000402964525
$ws.Range("C138").Value = 7026.178156858586
$ws.Range("C139").Value = 457.8330917196623
$ws.Range("C141").Value = 649.4459389945755

# --- Rows where column C was ".." (missing) -> now numeric 0 ---
$ws.Range("C100").Value = 0
$ws.Range("C119").Value = 0
$ws.Range("C140").Value = 0

# --- Update column AL (Colony) values 0 -> 1 ---
$ws.Range("AL7").Value = 1
$ws.Range("AL24").Value = 1
$ws.Range("AL26").Value = 1
$ws.Range("AL42").Value = 1
$ws.Range("AL44").Value = 1
$ws.Range("AL50").Value = 1
$ws.Range("AL51").Value = 1
$ws.Range("AL67").Value = 1
$ws.Range("AL68").Value = 1
$ws.Range("AL94").Value = 1
$ws.Range("AL98").Value = 1
$ws.Range("AL113").Value = 1
$ws.Range("AL117").Value = 1
$ws.Range("AL132").Value = 1
$ws.Range("AL138").Value = 1

Write-Output "edit applied"
